# edit.ps1 - apply the OOXML diff to before.docx using Word COM-interop style calls.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: in the paragraph "Nicolly-" change the trailing "-" run to ":"
# (keep it as its own run, just swap the text) and drop the _GoBack bookmark
# that currently sits there (it gets moved to the very end of the document).
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $ptext = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($ptext -eq "Nicolly-") {
        $dashRng = $p.Range
        $ok = $dashRng.Find.Execute("-", $true, $false, $false, $false, $false, `
                                     $true, 1, $false, "", 0)
        if ($ok) {
            $dashRng.Text = ":"
        }
    }
}

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# Change 2: the last (empty) paragraph gets "Falta:" typed into it, followed
# by three new paragraphs containing "2", "3" and "4". The _GoBack bookmark
# is re-created (collapsed) right after the "4".
# ---------------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)

# Build "2"
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$p2Index = $lastIndex + 1
$r2 = $d.Paragraphs($p2Index).Range
$r2.Collapse(0)
$r2.Text = "2"

# Build "3"
$r2b = $d.Paragraphs($p2Index).Range
$r2b.Collapse(0)
$r2b.InsertParagraphAfter()
$p3Index = $p2Index + 1
$r3 = $d.Paragraphs($p3Index).Range
$r3.Collapse(0)
$r3.Text = "3"

# Build "4" (with a temporary trailing placeholder char so that a collapsed
# bookmark can be anchored right after the "4" without snapping back to the
# start of the document).
$r3b = $d.Paragraphs($p3Index).Range
$r3b.Collapse(0)
$r3b.InsertParagraphAfter()
$p4Index = $p3Index + 1
$r4 = $d.Paragraphs($p4Index).Range
$r4.Collapse(0)
$r4.Text = "4ZZ"

$p4Range = $d.Paragraphs($p4Index).Range
$bmPos = $p4Range.Start + 1
$bmRng = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRng)

$p4Range2 = $d.Paragraphs($p4Index).Range
$trimRng = $d.Range($p4Range2.End - 3, $p4Range2.End - 1)
$trimRng.Text = ""

# Finally put "Falta:" into the original last paragraph, fixing up its run
# size so it matches the surrounding 28pt / 28pt(cs) formatting.
$rFalta = $lastPara.Range
$rFalta.Collapse(0)
$rFalta.Text = "Falta:"
$rFixup = $lastPara.Range
$rFixup.Font.Size = 14
$rFixup.Font.SizeBi = 14
